$d = $word.ActiveDocument

# 1) Programa - Portuguese paragraph: split the long numbered text into 6 runs
#    separated by manual line breaks.
$find1 = "1. Estrutura molecular e ligação química: Teoria de ligação de valência, estrutura de compostos com C, N, O; Relação entre estrutura e propriedades fisico-químicas2. Orbitais moleculares e as moléculas de O2 e N2: Limitações da teoria de ligação de valência, reatividade diferenciada de O2 e N2, relevância do O2 em sistemas biológicos, espécies reativas de oxigênio3. Ácidos, bases e a correlação com os ligantes dos metais em solução: Afinidade das bases por metais de transição, equilíbrio químico em sistemas biológicos4. Complexos metálicos - teoria do campo cristalino: Teoria do campo cristalino e os compostos octaédricos e tetraédricos; íons de metais de transição em sistemas biológicos5. Sistemas biológicos de transporte: Transporte de O2 em mamíferos, transferência de elétrons dependente de metaloproteínas;6. Processos catalíticos - ácido/base e oxido-redução em metaloproteínas: Proteínas contendo íon Zn2+, peroxidases, oxidases."
$repl1 = "1. Estrutura molecular e ligação química: Teoria de ligação de valência, estrutura de compostos com C, N, O; Relação entre estrutura e propriedades fisico-químicas^l2. Orbitais moleculares e as moléculas de O2 e N2: Limitações da teoria de ligação de valência, reatividade diferenciada de O2 e N2, relevância do O2 em sistemas biológicos, espécies reativas de oxigênio^l3. Ácidos, bases e a correlação com os ligantes dos metais em solução: Afinidade das bases por metais de transição, equilíbrio químico em sistemas biológicos^l4. Complexos metálicos - teoria do campo cristalino: Teoria do campo cristalino e os compostos octaédricos e tetraédricos; íons de metais de transição em sistemas biológicos^l5. Sistemas biológicos de transporte: Transporte de O2 em mamíferos, transferência de elétrons dependente de metaloproteínas;^l6. Processos catalíticos - ácido/base e oxido-redução em metaloproteínas: Proteínas contendo íon Zn2+, peroxidases, oxidases."
$d.Content.Find.Execute($find1, $true, $false, $false, $false, $false, $true, 1, $false, $repl1, 2) | Out-Null

# 2) Programa - English (italic) paragraph: same split pattern.
$find2 = "1. Chemical bonds and molecular structure: Valence bond theory, structure of compounds containing C, N and O, correlation of chemical structure with physical-chemical properties.2. Molecular orbitals and the O2 and N2 molecules: Limitations of the valence bond theory, varied reactivity of O2 and N2 molecules, relevance of O2 in biological systems, reactive oxygen species3. Acid/base as related to quelating agents: Bases and transition metallic ions, chemical equilibrium in biological systems4. Metallic ion complexes - crystal field theory: Crystal field theory describing octahedral and tetrahedral compounds, transition metal ions in biological systems5.Transport in biological systems: Oxygen transport and electron transfer mediated by metallo-proteins6. Acid/base and oxi-redox in metallo-proteins: Zn2+ proteins, peroxidases, oxidases"
$repl2 = "1. Chemical bonds and molecular structure: Valence bond theory, structure of compounds containing C, N and O, correlation of chemical structure with physical-chemical properties.^l2. Molecular orbitals and the O2 and N2 molecules: Limitations of the valence bond theory, varied reactivity of O2 and N2 molecules, relevance of O2 in biological systems, reactive oxygen species^l3. Acid/base as related to quelating agents: Bases and transition metallic ions, chemical equilibrium in biological systems^l4. Metallic ion complexes - crystal field theory: Crystal field theory describing octahedral and tetrahedral compounds, transition metal ions in biological systems^l5.Transport in biological systems: Oxygen transport and electron transfer mediated by metallo-proteins^l6. Acid/base and oxi-redox in metallo-proteins: Zn2+ proteins, peroxidases, oxidases"
$d.Content.Find.Execute($find2, $true, $false, $false, $false, $false, $true, 1, $false, $repl2, 2) | Out-Null

# 3) Avaliação - Critério run: split into 3 text segments, with a double break
#    before the final sentence.
$find3 = "A Nota final (NF) será calculada da seguinte maneira:NF = (P1 + 2*P2)/3Sendo que para P2 a matéria será cumulativa do semestre."
$repl3 = "A Nota final (NF) será calculada da seguinte maneira:^lNF = (P1 + 2*P2)/3^l^lSendo que para P2 a matéria será cumulativa do semestre."
$d.Content.Find.Execute($find3, $true, $false, $false, $false, $false, $true, 1, $false, $repl3, 2) | Out-Null

# 4) Avaliação - Norma de recuperação run: split into 2 text segments with a
#    double break between them.
$find4 = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$repl4 = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: ^l^lMR = (NF + PR)/2"
$d.Content.Find.Execute($find4, $true, $false, $false, $false, $false, $true, 1, $false, $repl4, 2) | Out-Null

# 5) Bibliografia run: split into 2 text segments with a single break.
$find5 = "1. Atkins e Jones, Princípios de Química, 5a edição, Bookman, 20112. Shiver e Atikins, Química Inorgânica, 4a edição, Bookman, 2008"
$repl5 = "1. Atkins e Jones, Princípios de Química, 5a edição, Bookman, 2011^l2. Shiver e Atikins, Química Inorgânica, 4a edição, Bookman, 2008"
$d.Content.Find.Execute($find5, $true, $false, $false, $false, $false, $true, 1, $false, $repl5, 2) | Out-Null
